$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1169.1628
$ws.Range("I129").Value = 505.25
$ws.Range("J129").Value = 1320.9143
$ws.Range("K129").Value = 1515.75
$ws.Range("L129").Value = 3962.7429
$ws.Range("M129").Value = 3484.25
$ws.Range("N129").Value = -13962.7429
$ws.Range("H137").Value = 1215.3889
$ws.Range("I137").Value = 1110.375
$ws.Range("J137").Value = 1425.4166
$ws.Range("K137").Value = 3331.125
$ws.Range("L137").Value = 4276.2498
$ws.Range("M137").Value = -781.125
$ws.Range("N137").Value = -9376.2498
$ws.Range("H138").Value = 3851.8684
$ws.Range("I138").Value = 2779.611
$ws.Range("J138").Value = 4816.9
$ws.Range("K138").Value = 8338.832999999999
$ws.Range("L138").Value = 14450.7
$ws.Range("M138").Value = -3198.832999999999
$ws.Range("N138").Value = -24730.7
$ws.Range("H141").Value = 5462.7666
$ws.Range("I141").Value = 2117.4285
$ws.Range("K141").Value = 6352.2855
$ws.Range("M141").Value = -1172.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9304.68
$ws.Range("I32").Value = 9545.815000000001
$ws.Range("J32").Value = 7737.3
$ws.Range("K32").Value = 9545.815000000001
$ws.Range("L32").Value = 7737.3
$ws.Range("M32").Value = -9258.815000000001
$ws.Range("N32").Value = -8311.299999999999
$ws.Range("H63").Value = 5584.1665
$ws.Range("I63").Value = 4876.25
$ws.Range("J63").Value = 7000
$ws.Range("K63").Value = 4876.25
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = -4190.25
$ws.Range("N63").Value = -8372
$ws.Range("H66").Value = 5584.1665
$ws.Range("I66").Value = 4876.25
$ws.Range("J66").Value = 7000
$ws.Range("K66").Value = 24381.25
$ws.Range("L66").Value = 35000
$ws.Range("M66").Value = -20949.25
$ws.Range("N66").Value = -41864
$ws.Range("H122").Value = 1765.2
$ws.Range("I122").Value = 1765.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5295.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2845.6
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5959.8887
$ws.Range("I22").Value = 6604.875
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 6604.875
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -6431.875
$ws.Range("N22").Value = -1146
$ws.Range("H43").Value = 130000
$ws.Range("J43").Value = 130000
$ws.Range("L43").Value = 130000
$ws.Range("N43").Value = -130362
$ws.Range("H134").Value = 5185.017
$ws.Range("I134").Value = 1885.2954
$ws.Range("J134").Value = 14864.2
$ws.Range("K134").Value = 5655.8862
$ws.Range("L134").Value = 44592.60000000001
$ws.Range("M134").Value = -3120.8862
$ws.Range("N134").Value = -49662.60000000001
$ws.Range("H138").Value = 58028.332
$ws.Range("J138").Value = 58028.332
$ws.Range("L138").Value = 58028.332
$ws.Range("N138").Value = -68308.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4814.8184
$ws.Range("J22").Value = 535.5
$ws.Range("L22").Value = 535.5
$ws.Range("N22").Value = -1235.5
$ws.Range("H31").Value = 3064.3333
$ws.Range("I31").Value = 2295.1875
$ws.Range("K31").Value = 2295.1875
$ws.Range("M31").Value = -2000.1875
$ws.Range("H34").Value = 3064.3333
$ws.Range("I34").Value = 2295.1875
$ws.Range("K34").Value = 2295.1875
$ws.Range("M34").Value = -2093.1875
$ws.Range("H94").Value = 10968.25
$ws.Range("I94").Value = 1723.8
$ws.Range("J94").Value = 15170.272
$ws.Range("K94").Value = 1723.8
$ws.Range("L94").Value = 15170.272
$ws.Range("M94").Value = -1272.8
$ws.Range("N94").Value = -16072.272
$ws.Range("H132").Value = 256225.84
$ws.Range("I132").Value = 322648.38
$ws.Range("J132").Value = 2612.5454
$ws.Range("K132").Value = 967945.14
$ws.Range("L132").Value = 7837.6362
$ws.Range("M132").Value = -965415.14
$ws.Range("N132").Value = -12897.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1498.8
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 1969.7142
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 5909.142599999999
$ws.Range("M32").Value = -917
$ws.Range("N32").Value = -6475.142599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 14442.714
$ws.Range("J51").Value = 14442.714
$ws.Range("L51").Value = 14442.714
$ws.Range("N51").Value = -15460.714
$ws.Range("H70").Value = 5375.5625
$ws.Range("I70").Value = 4817.278
$ws.Range("J70").Value = 6093.357
$ws.Range("K70").Value = 4817.278
$ws.Range("L70").Value = 6093.357
$ws.Range("M70").Value = -4547.278
$ws.Range("N70").Value = -6633.357
$ws.Range("H73").Value = 5375.5625
$ws.Range("I73").Value = 4817.278
$ws.Range("J73").Value = 6093.357
$ws.Range("K73").Value = 4817.278
$ws.Range("L73").Value = 6093.357
$ws.Range("M73").Value = -3881.278
$ws.Range("N73").Value = -7965.357
$ws.Range("H80").Value = 3033.7222
$ws.Range("I80").Value = 3180.5
$ws.Range("J80").Value = 2850.25
$ws.Range("K80").Value = 3180.5
$ws.Range("L80").Value = 2850.25
$ws.Range("M80").Value = -2182.5
$ws.Range("N80").Value = -4846.25
$ws.Range("H83").Value = 3033.7222
$ws.Range("I83").Value = 3180.5
$ws.Range("J83").Value = 2850.25
$ws.Range("K83").Value = 15902.5
$ws.Range("L83").Value = 14251.25
$ws.Range("M83").Value = -10910.5
$ws.Range("N83").Value = -24235.25
$ws.Range("H131").Value = 30866.166
$ws.Range("J131").Value = 30866.166
$ws.Range("L131").Value = 30866.166
$ws.Range("N131").Value = -40946.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 39599.5
$ws.Range("J117").Value = 39599.5
$ws.Range("L117").Value = 39599.5
$ws.Range("N117").Value = -48777.5
$ws.Range("H122").Value = 23689006
$ws.Range("I122").Value = 22732638
$ws.Range("J122").Value = 25004012
$ws.Range("K122").Value = 68197914
$ws.Range("L122").Value = 75012036
$ws.Range("M122").Value = -68195464
$ws.Range("N122").Value = -75016936
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23114
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 23114
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H118").Value = 30739
$ws.Range("J118").Value = 30739
$ws.Range("L118").Value = 30739
$ws.Range("N118").Value = -34053
$ws.Range("H123").Value = 24422.334
$ws.Range("J123").Value = 24422.334
$ws.Range("L123").Value = 24422.334
$ws.Range("N123").Value = -34222.334
$ws.Range("H136").Value = 874.0599999999999
$ws.Range("I136").Value = 968.2439000000001
$ws.Range("J136").Value = 445
$ws.Range("K136").Value = 2904.7317
$ws.Range("L136").Value = 1335
$ws.Range("M136").Value = -354.7317000000003
$ws.Range("N136").Value = -6435
